{"js": "// Replace the division-problem text in the worksheet table, cell by cell,\n// in document order. Each old value is unique in the document, so a\n// matchWholeWord search reliably finds the single run that needs updating.\n// Doing the replacements in this order also matters: the new value\n// \"145\u00f75=\" (introduced for the old \"792\u00f75=\" cell) must only be written\n// AFTER the original \"145\u00f75=\" cell has already been changed to \"267\u00f72=\",\n// otherwise the freshly-written text could be matched again.\nconst mapping = [\n  [\"573\u00f77=\", \"486\u00f73=\"],\n  [\"607\u00f72=\", \"274\u00f76=\"],\n  [\"337\u00f79=\", \"256\u00f76=\"],\n  [\"405\u00f72=\", \"655\u00f75=\"],\n  [\"295\u00f73=\", \"196\u00f77=\"],\n  [\"187\u00f78=\", \"598\u00f75=\"],\n  [\"854\u00f78=\", \"822\u00f78=\"],\n  [\"133\u00f76=\", \"820\u00f73=\"],\n  [\"186\u00f77=\", \"676\u00f72=\"],\n  [\"145\u00f75=\", \"267\u00f72=\"],\n  [\"210\u00f73=\", \"248\u00f74=\"],\n  [\"251\u00f76=\", \"584\u00f75=\"],\n  [\"706\u00f79=\", \"739\u00f72=\"],\n  [\"532\u00f72=\", \"496\u00f72=\"],\n  [\"998\u00f77=\", \"513\u00f79=\"],\n  [\"838\u00f77=\", \"689\u00f79=\"],\n  [\"841\u00f77=\", \"635\u00f72=\"],\n  [\"792\u00f75=\", \"145\u00f75=\"],\n  [\"230\u00f76=\", \"549\u00f75=\"],\n  [\"913\u00f78=\", \"391\u00f76=\"],\n  [\"108\u00f77=\", \"796\u00f76=\"],\n  [\"211\u00f74=\", \"162\u00f75=\"],\n  [\"689\u00f76=\", \"141\u00f76=\"],\n  [\"409\u00f75=\", \"457\u00f73=\"],\n  [\"955\u00f77=\", \"536\u00f73=\"],\n];\n\nfor (const [oldText, newText] of mapping) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: true,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  // Replace in place so the run's formatting (font, size, etc.) is kept.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Replace the division-problem text in the worksheet table, cell by cell,\n# in document order, using Find & Replace on the whole document range.\n# wdFindContinue = 1, wdReplaceOne = 2 (replace exactly the single match\n# located by this Execute call).\n#\n# The replacement order matters: the new value \"145\u00f75=\" (introduced for\n# the old \"792\u00f75=\" cell) must only be written AFTER the original \"145\u00f75=\"\n# cell has already been changed to \"267\u00f72=\", otherwise Find could match\n# the freshly-written text instead of the intended original cell.\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceOne = 2\n\n$pairs = @(\n    @(\"573\u00f77=\", \"486\u00f73=\"),\n    @(\"607\u00f72=\", \"274\u00f76=\"),\n    @(\"337\u00f79=\", \"256\u00f76=\"),\n    @(\"405\u00f72=\", \"655\u00f75=\"),\n    @(\"295\u00f73=\", \"196\u00f77=\"),\n    @(\"187\u00f78=\", \"598\u00f75=\"),\n    @(\"854\u00f78=\", \"822\u00f78=\"),\n    @(\"133\u00f76=\", \"820\u00f73=\"),\n    @(\"186\u00f77=\", \"676\u00f72=\"),\n    @(\"145\u00f75=\", \"267\u00f72=\"),\n    @(\"210\u00f73=\", \"248\u00f74=\"),\n    @(\"251\u00f76=\", \"584\u00f75=\"),\n    @(\"706\u00f79=\", \"739\u00f72=\"),\n    @(\"532\u00f72=\", \"496\u00f72=\"),\n    @(\"998\u00f77=\", \"513\u00f79=\"),\n    @(\"838\u00f77=\", \"689\u00f79=\"),\n    @(\"841\u00f77=\", \"635\u00f72=\"),\n    @(\"792\u00f75=\", \"145\u00f75=\"),\n    @(\"230\u00f76=\", \"549\u00f75=\"),\n    @(\"913\u00f78=\", \"391\u00f76=\"),\n    @(\"108\u00f77=\", \"796\u00f76=\"),\n    @(\"211\u00f74=\", \"162\u00f75=\"),\n    @(\"689\u00f76=\", \"141\u00f76=\"),\n    @(\"409\u00f75=\", \"457\u00f73=\"),\n    @(\"955\u00f77=\", \"536\u00f73=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $ok = $find.Execute($oldText, $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceOne)\n\n    if (-not $ok) {\n        Write-Output (\"NOT FOUND: \" + $oldText)\n    }\n}\n"}
